$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3119.9443
$ws.Range("J17").Value = 3119.9443
$ws.Range("L17").Value = 9359.832900000001
$ws.Range("N17").Value = -9695.832900000001
$ws.Range("H33").Value = 714755.3
$ws.Range("I33").Value = 833714.5
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 833714.5
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -833485.5
$ws.Range("N33").Value = -1458
$ws.Range("H53").Value = 860.5714
$ws.Range("I53").Value = 1449.5
$ws.Range("J53").Value = 418.875
$ws.Range("K53").Value = 1449.5
$ws.Range("L53").Value = 418.875
$ws.Range("M53").Value = -812.5
$ws.Range("N53").Value = -1692.875
$ws.Range("H64").Value = 90582670
$ws.Range("J64").Value = 4266
$ws.Range("L64").Value = 4266
$ws.Range("N64").Value = -4762
$ws.Range("H67").Value = 90582670
$ws.Range("J67").Value = 4266
$ws.Range("L67").Value = 4266
$ws.Range("N67").Value = -5982
$ws.Range("H87").Value = 82049.8
$ws.Range("J87").Value = 82049.8
$ws.Range("L87").Value = 82049.8
$ws.Range("N87").Value = -84545.8
$ws.Range("H90").Value = 82049.8
$ws.Range("J90").Value = 82049.8
$ws.Range("L90").Value = 246149.4
$ws.Range("N90").Value = -258629.4
$ws.Range("H100").Value = 1293.25
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 9394.687
$ws.Range("I132").Value = 4766.696
$ws.Range("K132").Value = 14300.088
$ws.Range("M132").Value = -11770.088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23966.104
$ws.Range("I32").Value = 27461.55
$ws.Range("K32").Value = 27461.55
$ws.Range("M32").Value = -27174.55
$ws.Range("H45").Value = 3396.2856
$ws.Range("I45").Value = 2727.7778
$ws.Range("J45").Value = 4599.6
$ws.Range("K45").Value = 2727.7778
$ws.Range("L45").Value = 4599.6
$ws.Range("M45").Value = -2350.7778
$ws.Range("N45").Value = -5353.6
$ws.Range("H80").Value = 77247
$ws.Range("J80").Value = 77247
$ws.Range("L80").Value = 77247
$ws.Range("N80").Value = -79243
$ws.Range("H83").Value = 77247
$ws.Range("J83").Value = 77247
$ws.Range("L83").Value = 231741
$ws.Range("N83").Value = -241725
$ws.Range("H132").Value = 14869.804
$ws.Range("I132").Value = 16293.711
$ws.Range("J132").Value = 10707.615
$ws.Range("K132").Value = 48881.133
$ws.Range("L132").Value = 32122.845
$ws.Range("M132").Value = -46351.133
$ws.Range("N132").Value = -37182.845
$ws.Range("H135").Value = 20000
$ws.Range("J135").Value = 20000
$ws.Range("L135").Value = 20000
$ws.Range("N135").Value = -30140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2535.8
$ws.Range("I86").Value = 2150.889
$ws.Range("K86").Value = 2150.889
$ws.Range("M86").Value = -1027.889
$ws.Range("H89").Value = 2535.8
$ws.Range("I89").Value = 2150.889
$ws.Range("K89").Value = 10754.445
$ws.Range("M89").Value = -5138.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1352.2727
$ws.Range("I16").Value = 1395.7646
$ws.Range("K16").Value = 1395.7646
$ws.Range("M16").Value = -1108.7646
$ws.Range("H31").Value = 16133454
$ws.Range("I31").Value = 31251966
$ws.Range("J31").Value = 7042.433
$ws.Range("K31").Value = 31251966
$ws.Range("L31").Value = 7042.433
$ws.Range("M31").Value = -31251671
$ws.Range("N31").Value = -7632.433
$ws.Range("H34").Value = 16133454
$ws.Range("I34").Value = 31251966
$ws.Range("J34").Value = 7042.433
$ws.Range("K34").Value = 31251966
$ws.Range("L34").Value = 7042.433
$ws.Range("M34").Value = -31251764
$ws.Range("N34").Value = -7446.433
$ws.Range("H74").Value = 63251
$ws.Range("I74").Value = 59000
$ws.Range("J74").Value = 64101.2
$ws.Range("K74").Value = 59000
$ws.Range("L74").Value = 64101.2
$ws.Range("M74").Value = -58126
$ws.Range("N74").Value = -65849.2
$ws.Range("H77").Value = 63251
$ws.Range("I77").Value = 59000
$ws.Range("J77").Value = 64101.2
$ws.Range("K77").Value = 177000
$ws.Range("L77").Value = 192303.6
$ws.Range("M77").Value = -172632
$ws.Range("N77").Value = -201039.6
$ws.Range("H103").Value = 31020
$ws.Range("I103").Value = 10682.667
$ws.Range("K103").Value = 10682.667
$ws.Range("M103").Value = -9510.666999999999
$ws.Range("H113").Value = 1352.2727
$ws.Range("I113").Value = 1395.7646
$ws.Range("K113").Value = 1395.7646
$ws.Range("M113").Value = 774.2354
$ws.Range("H120").Value = 45206.445
$ws.Range("J120").Value = 45206.445
$ws.Range("L120").Value = 45206.445
$ws.Range("N120").Value = -52464.445
$ws.Range("H124").Value = 44317
$ws.Range("J124").Value = 44317
$ws.Range("L124").Value = 44317
$ws.Range("N124").Value = -49227
$ws.Range("H132").Value = 33340964
$ws.Range("I132").Value = 37039530
$ws.Range("K132").Value = 111118590
$ws.Range("M132").Value = -111116060
$ws.Range("H134").Value = 1957.1578
$ws.Range("I134").Value = 1937.8064
$ws.Range("K134").Value = 5813.4192
$ws.Range("M134").Value = -3278.4192

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 8554.143
$ws.Range("J49").Value = 11760.4
$ws.Range("L49").Value = 35281.2
$ws.Range("N49").Value = -35593.2
$ws.Range("H107").Value = 7936704.5
$ws.Range("I107").Value = 394
$ws.Range("J107").Value = 15873015
$ws.Range("K107").Value = 1182
$ws.Range("L107").Value = 47619045
$ws.Range("M107").Value = 738
$ws.Range("N107").Value = -47622885
$ws.Range("H109").Value = 4517.4736
$ws.Range("I109").Value = 1694.7693
$ws.Range("J109").Value = 10633.333
$ws.Range("K109").Value = 5084.3079
$ws.Range("L109").Value = 31899.999
$ws.Range("M109").Value = -4044.3079
$ws.Range("N109").Value = -33979.999
$ws.Range("H112").Value = 8750
$ws.Range("I112").Value = 625
$ws.Range("K112").Value = 1875
$ws.Range("M112").Value = -767
$ws.Range("H122").Value = 1096.8572
$ws.Range("J122").Value = 1487.375
$ws.Range("L122").Value = 13386.375
$ws.Range("N122").Value = -18286.375
$ws.Range("H124").Value = 7734.5
$ws.Range("I124").Value = 472
$ws.Range("J124").Value = 14997
$ws.Range("K124").Value = 1416
$ws.Range("L124").Value = 44991
$ws.Range("M124").Value = 3494
$ws.Range("N124").Value = -54811
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = 1920
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 28500
$ws.Range("I126").Value = 28500
$ws.Range("K126").Value = 85500
$ws.Range("M126").Value = -80560
$ws.Range("H129").Value = 1194.8889
$ws.Range("I129").Value = 922.1429000000001
$ws.Range("J129").Value = 2149.5
$ws.Range("K129").Value = 2766.4287
$ws.Range("L129").Value = 6448.5
$ws.Range("M129").Value = 2233.5713
$ws.Range("N129").Value = -16448.5
$ws.Range("H131").Value = 19257392
$ws.Range("J131").Value = 21046624
$ws.Range("L131").Value = 63139872
$ws.Range("N131").Value = -63149952
$ws.Range("H132").Value = 4365.6665
$ws.Range("I132").Value = 1009.1
$ws.Range("J132").Value = 8561.375
$ws.Range("K132").Value = 9081.9
$ws.Range("L132").Value = 77052.375
$ws.Range("M132").Value = -6551.9
$ws.Range("N132").Value = -82112.375
$ws.Range("H133").Value = 25548.312
$ws.Range("J133").Value = 31645.834
$ws.Range("L133").Value = 94937.50199999999
$ws.Range("N133").Value = -105057.502

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4589.1562
$ws.Range("I7").Value = 3215.3333
$ws.Range("J7").Value = 5801.353
$ws.Range("K7").Value = 3215.3333
$ws.Range("L7").Value = 5801.353
$ws.Range("M7").Value = -3103.3333
$ws.Range("N7").Value = -6025.353
$ws.Range("H93").Value = 1966.6666
$ws.Range("I93").Value = 1966.6666
$ws.Range("K93").Value = 1966.6666
$ws.Range("M93").Value = -718.6666
$ws.Range("H122").Value = 16117.134
$ws.Range("I122").Value = 38000
$ws.Range("J122").Value = 12750.538
$ws.Range("K122").Value = 114000
$ws.Range("L122").Value = 38251.614
$ws.Range("M122").Value = -111550
$ws.Range("N122").Value = -43151.614
$ws.Range("H126").Value = 4589.1562
$ws.Range("I126").Value = 3215.3333
$ws.Range("J126").Value = 5801.353
$ws.Range("K126").Value = 9645.999899999999
$ws.Range("L126").Value = 17404.059
$ws.Range("M126").Value = -7175.999899999999
$ws.Range("N126").Value = -22344.059
$ws.Range("H132").Value = 6816.1816
$ws.Range("I132").Value = 6298.385
$ws.Range("J132").Value = 7564.1113
$ws.Range("K132").Value = 18895.155
$ws.Range("L132").Value = 22692.3339
$ws.Range("M132").Value = -16365.155
$ws.Range("N132").Value = -27752.3339
$ws.Range("H136").Value = 3150.0322
$ws.Range("I136").Value = 1686.7368
$ws.Range("K136").Value = 5060.2104
$ws.Range("M136").Value = -2510.2104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 100428.5
$ws.Range("J46").Value = 100428.5
$ws.Range("L46").Value = 100428.5
$ws.Range("N46").Value = -100890.5
$ws.Range("H107").Value = 537.6
$ws.Range("I107").Value = 547.5
$ws.Range("K107").Value = 1642.5
$ws.Range("M107").Value = 277.5
$ws.Range("H132").Value = 18521288
$ws.Range("I132").Value = 1329.25
$ws.Range("J132").Value = 45459410
$ws.Range("K132").Value = 3987.75
$ws.Range("L132").Value = 136378230
$ws.Range("M132").Value = -1457.75
$ws.Range("N132").Value = -136383290
$ws.Range("H134").Value = 100428.5
$ws.Range("J134").Value = 100428.5
$ws.Range("L134").Value = 301285.5
$ws.Range("N134").Value = -306355.5
